# Hortaliza, Femacal de La Calera - Acelga
# A new weekly price record was inserted as row 200, pushing the existing
# rows (200-247) down to (201-248).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 200, shifting rows 200:247 down to 201:248.
$ws.Rows("200:200").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(200, 1).Value  = 3
$ws.Cells.Item(200, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(200, 3).Value  = "Coquimbo"
$ws.Cells.Item(200, 4).Value  = 44543
$ws.Cells.Item(200, 5).Value  = 5
$ws.Cells.Item(200, 6).Value  = 100112009
$ws.Cells.Item(200, 7).Value  = "Acelga"
$ws.Cells.Item(200, 8).Value  = "Sin especificar"
$ws.Cells.Item(200, 9).Value  = "Primera"
$ws.Cells.Item(200, 10).Value = 230
$ws.Cells.Item(200, 11).Value = 2000
$ws.Cells.Item(200, 12).Value = 2200
$ws.Cells.Item(200, 13).Value = 2096
$ws.Cells.Item(200, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(200, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(200, 16).Value = 349
$ws.Cells.Item(200, 17).Value = 6
$ws.Cells.Item(200, 18).Value = "Hortaliza"

# Apply the same date number format used by the rest of column D (style
# index 2 in the original file maps to the "YYYY-MM-DD HH:MM:SS" numFmt).
$ws.Cells.Item(200, 4).NumberFormat = $ws.Cells.Item(201, 4).NumberFormat
